$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - headers (A1 stays the same, update B1/C1, add D1)
$ws.Range("B1").Value = "Tuotteen väri"
$ws.Range("C1").Value = "Malliston nimi"
$ws.Range("D1").Value = "poista"

# Row 2
$ws.Range("A2").Value = "hammer123"
$ws.Range("B2").Value = "musta"
$ws.Range("C2").Value = "deluxe"

# Row 3
$ws.Range("A3").Value = "helmet123"
$ws.Range("B3").Value = "sininen"
$ws.Range("C3").Value = "sale"
$ws.Range("D3").ClearContents()

# Row 4 (new row)
$ws.Range("A4").Value = "helmet123"
$ws.Range("D4").Value = "X"

# Row 5 (new row)
$ws.Range("A5").Value = "ski1"
$ws.Range("B5").Value = "valkoinen"
$ws.Range("C5").Value = "winter"

# Update selection to match the final active cell
$ws.Range("B5").Select()
